$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Row 12
$ws.Range("G12").Value = 1240524717.0500007
$ws.Range("J12").Value = 3588029419

# Row 13
$ws.Range("G13").Value = 319819483.18000001
$ws.Range("J13").Value = 956934340.60000002

# Row 14
$ws.Range("G14").Value = 34063116.800000042
$ws.Range("J14").Value = 146268235.09999999

# Row 16
$ws.Range("G16").Value = -60834434.380000003
$ws.Range("J16").Value = 193292161.30000001

# Row 18 - G18 becomes a formula (was a static value)
$ws.Range("G18").Formula = "=SUM(G12:G17)"

# Row 19
$ws.Range("G19").Value = -379300000.00000012
$ws.Range("J19").Value = 1105900000

# Row 21 - G21 becomes a formula (was a static value)
$ws.Range("G21").Formula = "=SUM(G18:G20)"

# Row 22
$ws.Range("G22").Value = -20015625
$ws.Range("J22").Value = 57815625

# Row 26
$ws.Range("G26").Value = 1029174575.116062
$ws.Range("J26").Value = 1031977291

$excel.CalculateFullRebuild()
$wb.Save()
